# Recoupmentplan vorlaeufig ausgefuellt. Erloesvorschau die Stueckzahl bei Anzu von der
# Gesamtsumme ausgeschlossen.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vorlage")

# --- Fix typo in shared string: add missing closing parenthesis ---
$ws.Range("A5").Value = "(eigene Rückstellungen)"

# --- New "Gesamtprojektkosten" / "Eigenanteil" rows (29/30) ---
$ws.Range("A29").Value = "Gesamtprojektkosten"
$ws.Range("B29").Value = 262342.48
$ws.Range("B29").Font.Size = 8

$ws.Range("A30").Value = "Eigenanteil"
$ws.Range("B30").Formula = "=0.05*B29"
$ws.Range("B30").Font.Size = 8

# Match the font size used for the existing figures in the "Beispielzahlen" block
$ws.Range("B26").Font.Size = 8
$ws.Range("B27").Font.Size = 8

# --- Update dependent figures on the "Rückstellungen" block ---
$ws.Range("C6").Value = 1

$ws.Range("D10").Formula = "=B29/2"
$ws.Range("D11").Value = 13117.12
$ws.Range("E11").Formula = "=D11/D10"
$ws.Range("E10").Formula = "=1-E11"

# --- Top summary row formulas ---
$ws.Range("D3").Formula = "=B9"
$ws.Range("F3").Formula = "=B9+D13"
$ws.Range("B4").Formula = "=B6"
$ws.Range("D4").Formula = "=D3+D10+D11"

# --- View/zoom adjustments ---
[void]$ws.Range("E24").Select()
$ws.Application.ActiveWindow.Zoom = 115

$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 14
